# Update "program" schedule sheet for week 6 content refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("program")
$ws.Activate()

# ---------------------------------------------------------------------
# Week 1 (rows 2-8): date shifts from 9-Mar-2022 (44629) to 16-Mar-2022 (44636)
# ---------------------------------------------------------------------
$ws.Range("E2:E8").Value = 44636

# ---------------------------------------------------------------------
# Week 2 (rows 9-13): date text shifts from "16-03-2022" to "23-03-2022"
# ---------------------------------------------------------------------
$ws.Range("E9:E13").Value = "23-03-2022"

# ---------------------------------------------------------------------
# Row 14 (Assignment header for week 3): chapter label + date text change
# ---------------------------------------------------------------------
$ws.Cells.Item(14,3).Value = "Assignment for week 2"
$ws.Cells.Item(14,5).Value = "30-03-2022"

# ---------------------------------------------------------------------
# Week 3 content rows 15-16: date text shifts to "30-03-2022"
# ---------------------------------------------------------------------
$ws.Range("E15:E16").Value = "30-03-2022"

# ---------------------------------------------------------------------
# Row 17 (Assignment header for week 4): chapter label change, and the
# date switches from a text shared-string to a literal serial number
# ---------------------------------------------------------------------
$ws.Cells.Item(17,3).Value = "Assignment  for week 3"
$ws.Cells.Item(17,5).Value = 44716

# ---------------------------------------------------------------------
# Week 4 content rows 18-20: date switches to literal serial 44716
# ---------------------------------------------------------------------
$ws.Range("E18:E20").Value = 44716

# ---------------------------------------------------------------------
# Row 21 used to be the "Assignment 3" header row; it is replaced with a
# normal content row duplicating row 22 ("Spatial Data Visualization").
# First copy formatting from row 22 so the red/bold "assignment" style
# is replaced with the regular content style, then set new values.
# ---------------------------------------------------------------------
$ws.Range("B22:F22").Copy()
$ws.Range("B21:F21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(21,1).Value = 21
$ws.Cells.Item(21,2).Value = "Spatial Data Visualization with R "
$ws.Cells.Item(21,3).Value = "Spatial Data Visualization"
$ws.Cells.Item(21,4).Value = 5
$ws.Cells.Item(21,5).Value = 44664
$ws.Cells.Item(21,6).Value = 120

# Row 22 keeps the same topic but moves to the new date / wording
$ws.Cells.Item(22,2).Value = "Spatial Data Visualization with R "
$ws.Cells.Item(22,3).Value = "Spatial Data Visualization"
$ws.Cells.Item(22,5).Value = 44664

# ---------------------------------------------------------------------
# Row 23 (Assignment header for week 5/6 boundary): chapter + date text
# ---------------------------------------------------------------------
$ws.Cells.Item(23,3).Value = "Assignment for week 5"
$ws.Cells.Item(23,5).Value = "20-04-2022"

# ---------------------------------------------------------------------
# Week 6 content rows 24-26: brand new topics replacing "Modelling"
# ---------------------------------------------------------------------
$ws.Cells.Item(24,2).Value = "Example 1 with spatial urban data"
$ws.Cells.Item(24,3).Value = "Spatial Data Visualization and Computation"
$ws.Cells.Item(24,5).Value = "20-04-2022"

$ws.Cells.Item(25,2).Value = "Example 2 with spatial urban data"
$ws.Cells.Item(25,3).Value = "Spatial Data Visualization and Computation"
$ws.Cells.Item(25,5).Value = "20-04-2022"

$ws.Cells.Item(26,2).Value = "Introduction to agent based model"
$ws.Cells.Item(26,3).Value = "Spatial Data Visualization and Computation"
$ws.Cells.Item(26,5).Value = "20-04-2022"

# Column C on rows 24-26 now wraps (longer shared label)
$ws.Range("C24:C26").WrapText = $true
$ws.Rows(25).RowHeight = 30
$ws.Rows(26).RowHeight = 30

# ---------------------------------------------------------------------
# Row 27 (Assignment header for week 7): date text change only
# ---------------------------------------------------------------------
$ws.Cells.Item(27,5).Value = "27-04-2022"

# ---------------------------------------------------------------------
# Row 28: new topic + date text change
# ---------------------------------------------------------------------
$ws.Cells.Item(28,2).Value = "Introduction to network modelling for urban computation"
$ws.Cells.Item(28,3).Value = "Network computation"
$ws.Cells.Item(28,5).Value = "27-04-2022"

# ---------------------------------------------------------------------
# Row 29 ("Network visualization and optimization" / "Extra curriculum")
# is removed entirely.
# ---------------------------------------------------------------------
$ws.Rows(29).Delete()

# ---------------------------------------------------------------------
# Add a running total of minutes spent.
# ---------------------------------------------------------------------
$ws.Range("G28").Formula = "=SUM(F2:F28)"

# ---------------------------------------------------------------------
# Restore sheet view: drop the frozen top-left scroll position and move
# the active selection to D25.
# ---------------------------------------------------------------------
$ws.Range("D25").Select()
